# The author renamed the "Codigo_Material" column header (cell J1) to
# "Codigo" and left the selection on J2 (just below the edited header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Codigo"

$ws.Range("J2").Select()
